$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(53).Insert()

$ws.Cells.Item(53, 1).Value = 6
$ws.Cells.Item(53, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(53, 3).Value = "Metropolitana"
$ws.Cells.Item(53, 4).Value = 44704
$ws.Cells.Item(53, 5).Value = 13
$ws.Cells.Item(53, 6).Value = 100114007
$ws.Cells.Item(53, 7).Value = "Jengibre"
$ws.Cells.Item(53, 8).Value = "Sin especificar"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 380
$ws.Cells.Item(53, 11).Value = 10000
$ws.Cells.Item(53, 12).Value = 10000
$ws.Cells.Item(53, 13).Value = 10000
$ws.Cells.Item(53, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(53, 15).Value = "Perú"
$ws.Cells.Item(53, 16).Value = 769
$ws.Cells.Item(53, 17).Value = 13
$ws.Cells.Item(53, 18).Value = "Hortaliza"
